# Auto-generated Excel COM-interop script
# Applies updated currentAveragePrice / LevePrice / LeveProfit figures
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$changes = @{
  "ALC" = @{
    101 = @{
      "H" = 365.66666
      "I" = 365.66666
      "K" = 1096.99998
      "M" = 525.0000199999999
    }
    132 = @{
      "H" = 16856.018
      "I" = 2465.6
      "K" = 7396.799999999999
      "M" = -4866.799999999999
    }
  }
  "ARM" = @{
    2 = @{
      "H" = 1514.2927
      "I" = 1586.6562
      "J" = 1257
      "K" = 1586.6562
      "L" = 1257
      "M" = -1473.6562
      "N" = -1483
    }
    9 = @{
      "H" = 30000
      "J" = 30000
      "L" = 30000
      "N" = -30340
    }
    20 = @{
      "H" = 30000
      "J" = 30000
      "L" = 30000
      "N" = -30540
    }
    45 = @{
      "H" = 65182.156
      "I" = 68413.72
      "J" = 7014
      "K" = 68413.72
      "L" = 7014
      "M" = -68036.72
      "N" = -7768
    }
    63 = @{
      "H" = 3606.3635
      "I" = 2458.75
      "J" = 6666.6665
      "K" = 2458.75
      "L" = 6666.6665
      "M" = -1772.75
      "N" = -8038.6665
    }
    66 = @{
      "H" = 3606.3635
      "I" = 2458.75
      "J" = 6666.6665
      "K" = 12293.75
      "L" = 33333.3325
      "M" = -8861.75
      "N" = -40197.3325
    }
    80 = @{
      "H" = 47499.5
      "J" = 47499.5
      "L" = 47499.5
      "N" = -49495.5
    }
    83 = @{
      "H" = 47499.5
      "J" = 47499.5
      "L" = 142498.5
      "N" = -152482.5
    }
    116 = @{
      "H" = 1514.2927
      "I" = 1586.6562
      "J" = 1257
      "K" = 1586.6562
      "L" = 1257
      "M" = 707.3438000000001
      "N" = -5845
    }
  }
  "BSM" = @{
    3 = @{
      "H" = 1514.2927
      "I" = 1586.6562
      "J" = 1257
      "K" = 1586.6562
      "L" = 1257
      "M" = -1472.6562
      "N" = -1485
    }
    35 = @{
      "H" = 29702.5
      "J" = 29702.5
      "L" = 29702.5
      "N" = -30322.5
    }
    82 = @{
      "H" = 18888.889
      "I" = 10000
      "K" = 10000
      "M" = -9617
    }
    85 = @{
      "H" = 18888.889
      "I" = 10000
      "K" = 10000
      "M" = -8674
    }
    119 = @{
      "H" = 42500
      "J" = 42500
      "L" = 42500
      "N" = -52176
    }
  }
  "CRP" = @{
    22 = @{
      "H" = 3296.1428
      "I" = 424.33334
      "J" = 5450
      "K" = 424.33334
      "L" = 5450
      "M" = -74.33334000000002
      "N" = -6150
    }
    31 = @{
      "H" = 5057703.5
      "I" = 4361.6665
      "J" = 5563037.5
      "K" = 4361.6665
      "L" = 5563037.5
      "M" = -4066.6665
      "N" = -5563627.5
    }
    34 = @{
      "H" = 5057703.5
      "I" = 4361.6665
      "J" = 5563037.5
      "K" = 4361.6665
      "L" = 5563037.5
      "M" = -4159.6665
      "N" = -5563441.5
    }
    123 = @{
      "H" = 25000
      "J" = 25000
      "L" = 25000
      "N" = -34800
    }
  }
  "CUL" = @{
    49 = @{
      "H" = 2566.6667
      "J" = 2566.6667
      "L" = 7700.000100000001
      "N" = -8012.000100000001
    }
    62 = @{
      "H" = 2500
      "I" = 2500
      "K" = 7500
      "M" = -6814
    }
    65 = @{
      "H" = 2500
      "I" = 2500
      "K" = 22500
      "M" = -19068
    }
    131 = @{
      "H" = 884.88776
      "I" = 483.75
      "J" = 920.54443
      "K" = 1451.25
      "L" = 2761.63329
      "M" = 3588.75
      "N" = -12841.63329
    }
  }
  "GSM" = @{
    126 = @{
      "H" = 5283
      "I" = 7581
      "J" = 2219
      "K" = 22743
      "L" = 6657
      "M" = -20273
      "N" = -11597
    }
  }
  "LTW" = @{
    22 = @{
      "H" = 972.375
      "I" = 1300
      "J" = 863.1667
      "K" = 1300
      "L" = 863.1667
      "M" = -1005
      "N" = -1453.1667
    }
    27 = @{
      "H" = 972.375
      "I" = 1300
      "J" = 863.1667
      "K" = 1300
      "L" = 863.1667
      "M" = -1193
      "N" = -1077.1667
    }
    46 = @{
      "H" = 4912.4165
      "I" = 860.8461
      "J" = 7202.4346
      "K" = 860.8461
      "L" = 7202.4346
      "M" = -672.8461
      "N" = -7578.4346
    }
    55 = @{
      "H" = 802.1818
      "I" = 788.38464
      "J" = 822.1111
      "K" = 788.38464
      "L" = 822.1111
      "M" = -615.38464
      "N" = -1168.1111
    }
    127 = @{
      "H" = 44245
      "J" = 44245
      "L" = 44245
      "N" = -54165
    }
    133 = @{
      "H" = 43313.25
      "J" = 43313.25
      "L" = 43313.25
      "N" = -48373.25
    }
    134 = @{
      "H" = 62475.668
      "J" = 62475.668
      "L" = 62475.668
      "N" = -72615.66800000001
    }
    135 = @{
      "H" = 37843.332
      "J" = 37843.332
      "L" = 37843.332
      "N" = -47983.332
    }
    136 = @{
      "H" = 2201.2
      "I" = 1645.9375
      "K" = 4937.8125
      "M" = -2387.8125
    }
    137 = @{
      "H" = 44662.5
      "J" = 44662.5
      "L" = 44662.5
      "N" = -54862.5
    }
  }
  "WVR" = @{
    62 = @{
      "H" = 2460
      "I" = 2400
      "J" = 2475
      "K" = 2400
      "L" = 2475
      "M" = -1776
      "N" = -3723
    }
    65 = @{
      "H" = 2460
      "I" = 2400
      "J" = 2475
      "K" = 12000
      "L" = 12375
      "M" = -8880
      "N" = -18615
    }
    81 = @{
      "H" = 1500.125
      "I" = 1500.125
      "J" = 0
      "K" = 3000.25
      "L" = 0
      "M" = -1939.25
      "N" = "DELETE"
    }
    84 = @{
      "H" = 1500.125
      "I" = 1500.125
      "J" = 0
      "K" = 15001.25
      "L" = 0
      "M" = -9697.25
      "N" = "DELETE"
    }
    125 = @{
      "H" = 36633.332
      "J" = 36633.332
      "L" = 36633.332
      "N" = -46473.332
    }
    126 = @{
      "H" = 981585.9399999999
      "I" = 1132403
      "J" = 1275
      "K" = 3397209
      "L" = 3825
      "M" = -3394739
      "N" = -8765
    }
    131 = @{
      "H" = 0
      "J" = 0
      "L" = 0
      "N" = "DELETE"
    }
    132 = @{
      "H" = 2572.2856
      "I" = 2074.3333
      "J" = 3468.6
      "K" = 6222.999899999999
      "L" = 10405.8
      "M" = -3692.999899999999
      "N" = -15465.8
    }
    136 = @{
      "H" = 1270.7587
      "I" = 1070.6818
      "K" = 3212.0454
      "M" = -662.0454
    }
  }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $sheetChanges = $changes[$sheetName]
    foreach ($rowNum in $sheetChanges.Keys) {
        $rowChanges = $sheetChanges[$rowNum]
        foreach ($col in $rowChanges.Keys) {
            $addr = "$col$rowNum"
            $val = $rowChanges[$col]
            if ($val -eq "DELETE") {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}
